# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Albahaca" ahead of the existing row 31, shifting all subsequent rows down
# by one (dimension grows from A1:R103 to A1:R104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 31..103 down to 32..104, carrying their values/formatting along.
$ws.Rows.Item(31).Insert()

# Populate the newly-opened row 31 with the new record.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44624
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 100112052
$ws.Range("G31").Value = "Albahaca"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 140
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 5500
$ws.Range("N31").Value = "$/docena de matas"
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 917
$ws.Range("Q31").Value = 6
$ws.Range("R31").Value = "Hortaliza"
